# Conferences.xlsx edit:
#  - add a new summary row (29) that averages the x-coordinate column
#  - re-fit the rows whose height had been stretched by the superscript
#    ("th"/"st"/"nd") runs now that the runs no longer carry a charset
#    override (the rows simply return to the sheet's default height)
#  - leave the final selection on F20, matching the saved workbook state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row with the average of column A (rows 2 through 28)
$ws.Range("A29").Formula = "=AVERAGE(A2:A28)"

# Rows 4-19 had an inflated row height (13.4) because of superscript
# ("th"/"st"/"nd") runs; auto-fit puts them back to the default (12.8)
$ws.Range("A4:A19").EntireRow.AutoFit()

# Final selection/active cell used when the workbook was last saved
$ws.Range("F20").Select()
